# fix: Corrigir schema PostgreSQL para compatibilidade
# Adds three new lead rows (registros 17, 18, 19) to the sheet, growing
# the used range from A1:V17 to A1:V20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Several source columns hold values that look like dates or numbers
    # (e.g. "2025-07-14", "+551177778888", "7632684980") but must be
    # stored as plain text, exactly like the rest of the sheet. Briefly
    # mark the cell as Text before assigning the value so Excel does not
    # auto-convert it to a date/number, then restore the default style
    # so no visible formatting change is introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 18 (registro 17): Maria Silva / BUEIRO na Pracinha
# ---------------------------------------------------------------------
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 17
Set-TextValue $ws.Range("C18") "2025-07-14"
$ws.Range("E18").Value = "Maria Silva"
Set-TextValue $ws.Range("F18") "+551177778888"
$ws.Range("H18").Value = "Pracinha, esquina da Rua José da Silva Calvo"
$ws.Range("I18").Value = "BUEIRO"
$ws.Range("J18").Value = "Bueiro entupido na esquina da Rua José da Silva Calvo com cheiro de esgoto"
$ws.Range("K18").Value = "ALTA"
$ws.Range("L18").Value = $true
$ws.Range("M18").Value = "texto_digitado"
$ws.Range("N18").Value = 0.708
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 45858.74813476107
$ws.Range("Q18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("R18").Value = 45858.74811244453
$ws.Range("R18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
Set-TextValue $ws.Range("S18") "7632684980"
$ws.Range("T18").Value = "Hoje 14/07 falei com Maria Silva na Pracinha, bueiro entupido na esquina da Rua jose da silva calvo, telefone 11 7777-8888, urgente porque tem cheiro de esgoto"

# ---------------------------------------------------------------------
# Row 19 (registro 18): entrada de sistema contendo a connection string
# do PostgreSQL (o motivo deste commit)
# ---------------------------------------------------------------------
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 18
$ws.Range("M19").Value = "texto_digitado"
$ws.Range("N19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 45858.75657905206
$ws.Range("Q19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("R19").Value = 45858.7565560986
$ws.Range("R19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
Set-TextValue $ws.Range("S19") "7632684980"
$ws.Range("T19").Value = "postgresql://postgres:HMjrapQYaqXpKVYwrHbNtzRHHIdQhmJZ@shinkansen.proxy.rlwy.net:15314/railway"

# ---------------------------------------------------------------------
# Row 20 (registro 19): Tereza / FIACAO na Consolação
# ---------------------------------------------------------------------
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 19
Set-TextValue $ws.Range("C20") "2025-12-08"
$ws.Range("E20").Value = "Tereza"
Set-TextValue $ws.Range("F20") "+55119999977333"
$ws.Range("G20").Value = "Consolação"
$ws.Range("I20").Value = "FIACAO"
$ws.Range("J20").Value = "Problema com a fiação, risco de choque"
$ws.Range("K20").Value = "ALTA"
$ws.Range("M20").Value = "texto_digitado"
$ws.Range("N20").Value = 0.92
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 45858.75687256754
$ws.Range("Q20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("R20").Value = 45858.7568501126
$ws.Range("R20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
Set-TextValue $ws.Range("S20") "7632684980"
$ws.Range("T20").Value = "08/12 falei com  Tereza na consolacao, problema com a fiacao telefone 11 99999-77333, risco de choque"

Write-Host "Added rows 18-20 (registros 17-19) to $($ws.Name)"
